$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_detection_template")

# --- Row 3: standalone formulas, wrapped in TRIM(...) ---
$ws.Range("N3").Formula = '=TRIM(IF($D3<>"","Mean ","Any ")&IF($H3="",$G3,"")' + "`n" + '&IF($H3<>"",$H3,"")' + "`n" + '&IF(AND($L3="LT",I3="ST")," LT-neg ST-pos",' + "`n" + 'IF($I3<>""," "&$I3,"")' + "`n" + '&IF(OR($I3="LT",$I3="ST",$I3<>""),"-pos","")' + "`n" + '&IF($K3<>""," "&$K3,"")' + "`n" + '&IF($J3<>""," "&$J3&"-pos","")' + "`n" + '&IF($L3<>""," "&$L3&"-neg",""))&IF($D3="",""," "&$D3)&", by "&$C3)'

$ws.Range("O3").Formula = '=TRIM(IF($I3="",IF($H3="",$G3,$H3),$I3)&" aggregate data")'

$ws.Range("P3").Formula = '=TRIM(IF($E3="Eukaryota","Eukaryote",$E3)&" in "&$B3&" detection aggregate data")'

$ws.Range("Q3").Formula = '=TRIM("Aggregate organism in "&$B3&" detection data")'

# --- Rows 4:13 : set the same relative formula on the whole range at once so
#     the engine records them as a shared-formula group (matches the target
#     diff, which has N4:N13 / O4:O13 / P4:P13 / Q4:Q13 as shared groups). ---

$ws.Range("N4:N13").Formula = '=TRIM(IF($D4<>"","Mean ","Any ")&IF($H4="",$G4,"")' + "`n" + '&IF($H4<>"",$H4,"")' + "`n" + '&IF(AND($L4="LT",I4="ST")," LT-neg ST-pos",' + "`n" + 'IF($I4<>""," "&$I4,"")' + "`n" + '&IF(OR($I4="LT",$I4="ST",$I4<>""),"-pos","")' + "`n" + '&IF($K4<>""," "&$K4,"")' + "`n" + '&IF($J4<>""," "&$J4&"-pos","")' + "`n" + '&IF($L4<>""," "&$L4&"-neg",""))&IF($D4="",""," "&$D4)&", by "&$C4)'

$ws.Range("O4:O13").Formula = '=TRIM(IF($I4="",IF($H4="",$G4,$H4),$I4)&" aggregate data")'

$ws.Range("P4:P13").Formula = '=TRIM(IF($E4="Eukaryota","Eukaryote",$E4)&" in "&$B4&" detection aggregate data")'

$ws.Range("Q4:Q13").Formula = '=TRIM("Aggregate organism in "&$B4&" detection data")'

# --- Update frozen-pane view / selection ---
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("G8").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("R6").Select()

$wb.Save()
